$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Rabanito" that sits (by date)
# just before the current row 198, so insert a new row there and shift the
# rest of the table (198-278) down to (199-279).
$ws.Rows(198).Insert()

# Populate the newly inserted row with the new record. Columns that stay
# constant for every row of this sub-sheet are copied from the neighbouring
# rows; only the date (D) and volume (J) are genuinely new data points.
$ws.Cells.Item(198, 1).Value = 9
$ws.Cells.Item(198, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44726
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = 300000001
$ws.Cells.Item(198, 7).Value = "Rabanito"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 9000
$ws.Cells.Item(198, 11).Value = 3000
$ws.Cells.Item(198, 12).Value = 3000
$ws.Cells.Item(198, 13).Value = 3000
$ws.Cells.Item(198, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(198, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(198, 16).Value = 30
$ws.Cells.Item(198, 17).Value = 100
$ws.Cells.Item(198, 18).Value = "Hortaliza"
